$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern used below for every text-valued cell: write the literal
# text via a formula ("=""01""") so Excel's "looks like a number/percent"
# auto-coercion never kicks in, then Copy / PasteSpecial(xlPasteValues)
# to collapse the formula down to a plain shared-string text value. This
# avoids touching NumberFormat/Style (which would stamp a new, unwanted
# cell style onto the sheet).

# --- "Week" column (B) -----------------------------------------------
$ws.Range("B2").Formula = '="01"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

$ws.Range("B3").Formula = '="03"'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)

$ws.Range("B4").Formula = '="02"'
$ws.Range("B4").Copy()
$ws.Range("B4").PasteSpecial(-4163)

# --- Numeric columns (A, C, D, E) --------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3

$ws.Range("C2").Value = 94
$ws.Range("D2").Value = 37
$ws.Range("E2").Value = 57

$ws.Range("C3").Value = 75
$ws.Range("D3").Value = 24
$ws.Range("E3").Value = 51

$ws.Range("C4").Value = 85
$ws.Range("D4").Value = 30
$ws.Range("E4").Value = 55

# --- "Percent Correct" column (F) --------------------------------------
$ws.Range("F2").Formula = '="100%"'
$ws.Range("F2").Copy()
$ws.Range("F2").PasteSpecial(-4163)

$ws.Range("F3").Formula = '="100%"'
$ws.Range("F3").Copy()
$ws.Range("F3").PasteSpecial(-4163)

$ws.Range("F4").Formula = '="100%"'
$ws.Range("F4").Copy()
$ws.Range("F4").PasteSpecial(-4163)

$excel.CutCopyMode = $false
